$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("CalculatorTestsWithURI")

# New "Serialized" test-case rows (TC_Sub_Object / TC_Mul_Object / TC_Div_Object),
# mirroring the existing TC_Add_Object row (row 6) but with the headers/expected
# response body reused from the Subtract/Multiply/Divide rows (3/4/5) and an
# empty PostBody column (serialized-object tests have no PostBody).
$newRows = @(
    @{ row = 7; id = "TC_Sub_Object"; srcRow = 3 },
    @{ row = 8; id = "TC_Mul_Object"; srcRow = 4 },
    @{ row = 9; id = "TC_Div_Object"; srcRow = 5 }
)

foreach ($r in $newRows) {
    $rowNum = $r.row
    $srcRow = $r.srcRow

    # Seed the new row with row 6's formatting (same column styles / hyperlink cell)
    $ws.Range("A6:F6").Copy()
    $ws.Range("A" + $rowNum + ":F" + $rowNum).PasteSpecial(-4122)

    $ws.Range("A$rowNum").Value = $r.id
    $ws.Range("B$rowNum").Value = $ws.Range("B6").Value()
    $ws.Range("C$rowNum").Value = $ws.Range("C6").Value()
    $ws.Range("D$rowNum").Value = $ws.Range("D$srcRow").Value()
    $ws.Range("F$rowNum").Value = $ws.Range("F$srcRow").Value()
    # E (PostBody) intentionally left blank

    $ws.Hyperlinks.Add($ws.Range("B$rowNum"), "http://www.dneonline.com/calculator.asmx")
    # Re-apply the hyperlink-cell formatting (Hyperlinks.Add resets it)
    $ws.Range("B6").Copy()
    $ws.Range("B$rowNum").PasteSpecial(-4122)

    # Undo the auto row-height bump caused by the multi-line Headers/Body text
    $ws.Rows($rowNum).AutoFit()
}

# Grow the table / autofilter to cover the new rows
$lo = $ws.ListObjects.Item(1)
$lo.Resize($ws.Range("A1:F9"))

$ws.Range("D14").Select()
